$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of data is inserted above the existing rows (weekly refresh):
# insert a blank row at row 79, which shifts the old rows 79-157 down to 80-158,
# then populate the newly inserted row 79 with this week's record.
$ws.Rows(79).Insert()

$ws.Cells.Item(79, 1).Value = 8
$ws.Cells.Item(79, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(79, 3).Value = "Coquimbo"
$ws.Cells.Item(79, 4).Value = 45128
$ws.Cells.Item(79, 5).Value = 4
$ws.Cells.Item(79, 6).Value = 100114007
$ws.Cells.Item(79, 7).Value = "Jengibre"
$ws.Cells.Item(79, 8).Value = "Sin especificar"
$ws.Cells.Item(79, 9).Value = "Primera"
$ws.Cells.Item(79, 10).Value = 360
$ws.Cells.Item(79, 11).Value = 17000
$ws.Cells.Item(79, 12).Value = 17500
$ws.Cells.Item(79, 13).Value = 17250
$ws.Cells.Item(79, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(79, 15).Value = "Perú"
$ws.Cells.Item(79, 16).Value = 1327
$ws.Cells.Item(79, 17).Value = 13
$ws.Cells.Item(79, 18).Value = "Hortaliza"
